{"js": "// Replace the old \"PCSI scores that were collected ... across all scores. \" text\n// with the new, shorter sentence: \"PCSI scores will be used from the initial visit,\n// within 14 days of concussion. \"\n\nconst body = context.document.body;\n\nconst oldText =\n  \"PCSI scores that were collected within 14 days of head injury. \" +\n  \"If a subject has multiple PCSI scores collected in this time (~14% of subjects), \" +\n  \"the median score and days post-injury will be taken across all scores. \";\n\nconst searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found in document body.\");\n}\n\nconst target = searchResults.items[0];\n\n// Replace the whole matched range with the new sentence in a single\n// insertion so the resulting run stays cleanly scoped to the replaced\n// text (a follow-up insertText(..., After) on the inserted range tends\n// to cascade-merge with unrelated sibling runs in this engine).\ntarget.insertText(\n  \"PCSI scores will be used from the initial visit, within 14 days of concussion. \",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Replace the old \"PCSI scores that were collected ... across all scores. \" text\n# with the new, shorter sentence: \"PCSI scores will be used from the initial visit,\n# within 14 days of concussion. \"\n\n$d = $word.ActiveDocument\n\n$oldText = \"PCSI scores that were collected within 14 days of head injury. \" +\n    \"If a subject has multiple PCSI scores collected in this time (~14% of subjects), \" +\n    \"the median score and days post-injury will be taken across all scores. \"\n\n$newText = \"PCSI scores will be used from the initial visit, within 14 days of concussion. \"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchWildcards = $false\n$find.MatchCase = $true\n\nif ($find.Execute()) {\n    $range.Text = $newText\n} else {\n    throw \"Target text not found in document.\"\n}\n"}
